# Natmi following Dr Hou advice
# Ccl20-Cxcr3 LR-pair sheet: add "Neutro" sending cluster and split target
# cluster into M1 / M2, recomputing all derived statistics (4 data rows
# instead of 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous single data row (and any stale cells beyond it) before
# writing the new 4-row table so leftover values don't linger.
$ws.Range("A2:T5").ClearContents()

$headerCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
  @("FAPs",  "Ccl20", "Cxcr3", "M1", 3, 1, 1.019213,              3.057639,            0.928181872270981,  0.928181872270981,
    2, 0.6666666666666666, 0.7166990000000001, 2.150097,          0.2276207788704612,  0.2276207788704611,
    0.7304689378870001,   6.574220440983001,   0.2112734806997636, 0.2112734806997636),

  @("FAPs",  "Ccl20", "Cxcr3", "M2", 3, 1, 1.019213,              3.057639,            0.928181872270981,  0.928181872270981,
    3, 1,                 2.431954666666666,  7.295864,           0.7723792211295388, 0.7723792211295388,
    2.478679811677333,    22.308118305096,    0.7169083915712174, 0.7169083915712174),

  @("Neutro","Ccl20", "Cxcr3", "M1", 1, 0.3333333333333333, 0.07886166666666666, 0.236585,     0.07181812772901904, 0.07181812772901904,
    2, 0.6666666666666666, 0.7166990000000001, 2.150097,          0.2276207788704612,  0.2276207788704611,
    0.05652007763833334,  0.508680698745,     0.01634729817069758, 0.01634729817069757),

  @("Neutro","Ccl20", "Cxcr3", "M2", 1, 0.3333333333333333, 0.07886166666666666, 0.236585,     0.07181812772901904, 0.07181812772901904,
    3, 1,                 2.431954666666666,  7.295864,           0.7723792211295388, 0.7723792211295388,
    0.1917879982711111,   1.72609198444,      0.05547082955832146, 0.05547082955832146)
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $rowNum = $r + 2
  $rowValues = $data[$r]
  for ($c = 0; $c -lt $rowValues.Length; $c++) {
    $cellRef = "{0}{1}" -f $headerCols[$c], $rowNum
    $ws.Range($cellRef).Value = $rowValues[$c]
  }
}
